# Daily update at 8 AM UTC
# Appends the next day's row of data to the "Wins Over Time" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 97

# Match the date-serial number format used by the rest of column A
# (copy it from the previous row so the new cell gets the same style).
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

$ws.Cells.Item($newRow, 1).Value = 46046
$ws.Cells.Item($newRow, 2).Value = 225
$ws.Cells.Item($newRow, 3).Value = 230
$ws.Cells.Item($newRow, 4).Value = 222
